$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 308-313 with revised figures ---
$ws.Range("B308").Value = 184
$ws.Range("C308").Value = 840
$ws.Range("D308").Value = 27
$ws.Range("E308").Value = 187

$ws.Range("B309").Value = 153
$ws.Range("C309").Value = 692
$ws.Range("D309").Value = 21
$ws.Range("E309").Value = 164

$ws.Range("B310").Value = 157
$ws.Range("C310").Value = 707
$ws.Range("D310").Value = 24
$ws.Range("E310").Value = 165

$ws.Range("B311").Value = 334
$ws.Range("C311").Value = 1600
$ws.Range("D311").Value = 40
$ws.Range("E311").Value = 319

$ws.Range("C312").Value = 1188
$ws.Range("D312").Value = 7
$ws.Range("E312").Value = 174

$ws.Range("B313").Value = 97
$ws.Range("C313").Value = 1005
$ws.Range("D313").Value = 14
$ws.Range("E313").Value = 153

# --- Append new row 314 with the new scraped data point ---
# Force column A to be stored as text so the date-like string "11.01.2021"
# isn't auto-converted into a date serial number, matching the existing
# date column's inline/shared-string representation.
$ws.Range("A314").NumberFormat = "@"
$ws.Range("A314").Value = "11.01.2021"
# Reset the cell style to match the surrounding (unstyled) data cells.
$ws.Range("A314").Style = $ws.Range("A313").Style

$ws.Range("B314").Value = 41
$ws.Range("C314").Value = 938
$ws.Range("D314").Value = 13
$ws.Range("E314").Value = 56
